# Generate Report for Handback
# Regenerates the handback-status workbook: the first tracked file is
# re-processed (new id/hash/timestamps) and a second file is handed back
# and appended as a new row on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$OLD_GUID = "d963e7ae-6d89-4111-ae40-56042b4814f9"
$NEW_GUID1 = "859e9e9f-dcc1-4d34-9199-9ef345bb5a9b"
$NEW_GUID2 = "c70a462c-d133-44c2-9e86-4df3f6cc1309"

$NEW_HASH1 = "7e7dfbe0f84d2017e4fc0d48a5f5c419ef44971a"
$NEW_HASH2 = "a669482bb75e2c896168956a5a11bb1c33b09a5e"

$OVERVIEW_DATE = "2016-08-16 11:02:13"
$ZHCN_HANDOFF_DATE = "2016-08-16 11:02:05"
$ZHCN_HANDBACK_DATE = "2016-08-16 11:02:31"
$DEDE_HANDOFF_DATE = "2016-08-16 11:02:13"
$DEDE_HANDBACK_DATE = "2016-08-16 11:02:38"

$xlPasteValues = -4163

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

# Row 2: the first file got a new id + new "Latest HO Xliff Generate Date"
$wsO.Range("A2").Value = "$NEW_GUID1.md"
$wsO.Range("B2").Value = "e2e\$NEW_GUID1.md"
$wsO.Range("G2").Value = $OVERVIEW_DATE

$wsO.Range("B2").Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cad10e00046eb8b3f7a6b0b07ba478fbb0ae7a70/e2e/$NEW_GUID1.md", "", "", "e2e\$NEW_GUID1.md") | Out-Null

# Row 3: second file, newly handed back
$wsO.Range("A3").Value = "$NEW_GUID2.md"
$wsO.Range("B3").Value = "e2e\$NEW_GUID2.md"
$wsO.Range("C2").Copy()
$wsO.Range("C3").PasteSpecial($xlPasteValues)
$wsO.Range("E2").Copy()
$wsO.Range("E3").PasteSpecial($xlPasteValues)
$wsO.Range("F2").Copy()
$wsO.Range("F3").PasteSpecial($xlPasteValues)
$wsO.Range("G3").Value = $OVERVIEW_DATE

$wsO.Range("B3").Style = "HyperLink"
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cad10e00046eb8b3f7a6b0b07ba478fbb0ae7a70/e2e/$NEW_GUID2.md", "", "", "e2e\$NEW_GUID2.md") | Out-Null

$wsO.ListObjects.Item(1).Resize($wsO.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 2 updates (same file, re-processed: new id/hash/timestamps)
$ws2.Range("A2").Value = "$NEW_GUID1.md"
$ws2.Range("G2").Value = "$NEW_GUID1.$NEW_HASH1.zh-cn.xlf"
$ws2.Range("H2").Value = $ZHCN_HANDOFF_DATE
$ws2.Range("I2").Value = "$NEW_GUID1.md"
$ws2.Range("J2").Value = "$NEW_GUID1.$NEW_HASH1.zh-cn.xlf"
$ws2.Range("K2").Value = $ZHCN_HANDBACK_DATE

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cad10e00046eb8b3f7a6b0b07ba478fbb0ae7a70/e2e/$NEW_GUID1.md", "", "", "$NEW_GUID1.md") | Out-Null
$ws2.Range("I2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/af7aca506d4b4407638c287aed10c67349926c33/e2e/$NEW_GUID1.md", "", "", "$NEW_GUID1.md") | Out-Null

# Row 3: new file
$ws2.Range("A3").Value = "$NEW_GUID2.md"
$ws2.Range("B2").Copy()
$ws2.Range("B3").PasteSpecial($xlPasteValues)
$ws2.Range("C2").Copy()
$ws2.Range("C3").PasteSpecial($xlPasteValues)
$ws2.Range("D2").Copy()
$ws2.Range("D3").PasteSpecial($xlPasteValues)
$ws2.Range("E2").Copy()
$ws2.Range("E3").PasteSpecial($xlPasteValues)
$ws2.Range("M2").Copy()
$ws2.Range("F3").PasteSpecial($xlPasteValues)
$ws2.Range("G3").Value = "$NEW_GUID2.$NEW_HASH2.zh-cn.xlf"
$ws2.Range("H3").Value = $ZHCN_HANDOFF_DATE
$ws2.Range("I3").Value = "$NEW_GUID2.md"
$ws2.Range("J3").Value = "$NEW_GUID2.$NEW_HASH2.zh-cn.xlf"
$ws2.Range("K3").Value = $ZHCN_HANDBACK_DATE
$ws2.Range("L2").Copy()
$ws2.Range("L3").PasteSpecial($xlPasteValues)
$ws2.Range("M2").Copy()
$ws2.Range("M3").PasteSpecial($xlPasteValues)
$ws2.Range("L2").Copy()
$ws2.Range("N3").PasteSpecial($xlPasteValues)
$ws2.Range("F2").Copy()
$ws2.Range("O3").PasteSpecial($xlPasteValues)
$ws2.Range("L2").Copy()
$ws2.Range("P3").PasteSpecial($xlPasteValues)

$ws2.Range("A3").Style = "HyperLink"
$ws2.Range("I3").Style = "HyperLink"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cad10e00046eb8b3f7a6b0b07ba478fbb0ae7a70/e2e/$NEW_GUID2.md", "", "", "$NEW_GUID2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/af7aca506d4b4407638c287aed10c67349926c33/e2e/$NEW_GUID2.md", "", "", "$NEW_GUID2.md") | Out-Null

$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 2 updates (same file, re-processed: new id/hash/timestamps)
$ws3.Range("A2").Value = "$NEW_GUID1.md"
$ws3.Range("G2").Value = "$NEW_GUID1.$NEW_HASH1.de-de.xlf"
$ws3.Range("H2").Value = $DEDE_HANDOFF_DATE
$ws3.Range("I2").Value = "$NEW_GUID1.md"
$ws3.Range("J2").Value = "$NEW_GUID1.$NEW_HASH1.de-de.xlf"
$ws3.Range("K2").Value = $DEDE_HANDBACK_DATE

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cad10e00046eb8b3f7a6b0b07ba478fbb0ae7a70/e2e/$NEW_GUID1.md", "", "", "$NEW_GUID1.md") | Out-Null
$ws3.Range("I2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/08fa60f54db3d798352b67ce07671559d3c9528e/e2e/$NEW_GUID1.md", "", "", "$NEW_GUID1.md") | Out-Null

# Row 3: new file
$ws3.Range("A3").Value = "$NEW_GUID2.md"
$ws3.Range("B2").Copy()
$ws3.Range("B3").PasteSpecial($xlPasteValues)
$ws3.Range("C2").Copy()
$ws3.Range("C3").PasteSpecial($xlPasteValues)
$ws3.Range("D2").Copy()
$ws3.Range("D3").PasteSpecial($xlPasteValues)
$ws3.Range("E2").Copy()
$ws3.Range("E3").PasteSpecial($xlPasteValues)
$ws3.Range("M2").Copy()
$ws3.Range("F3").PasteSpecial($xlPasteValues)
$ws3.Range("G3").Value = "$NEW_GUID2.$NEW_HASH2.de-de.xlf"
$ws3.Range("H3").Value = $DEDE_HANDOFF_DATE
$ws3.Range("I3").Value = "$NEW_GUID2.md"
$ws3.Range("J3").Value = "$NEW_GUID2.$NEW_HASH2.de-de.xlf"
$ws3.Range("K3").Value = $DEDE_HANDBACK_DATE
$ws3.Range("L2").Copy()
$ws3.Range("L3").PasteSpecial($xlPasteValues)
$ws3.Range("M2").Copy()
$ws3.Range("M3").PasteSpecial($xlPasteValues)
$ws3.Range("L2").Copy()
$ws3.Range("N3").PasteSpecial($xlPasteValues)
$ws3.Range("F2").Copy()
$ws3.Range("O3").PasteSpecial($xlPasteValues)
$ws3.Range("L2").Copy()
$ws3.Range("P3").PasteSpecial($xlPasteValues)

$ws3.Range("A3").Style = "HyperLink"
$ws3.Range("I3").Style = "HyperLink"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cad10e00046eb8b3f7a6b0b07ba478fbb0ae7a70/e2e/$NEW_GUID2.md", "", "", "$NEW_GUID2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/08fa60f54db3d798352b67ce07671559d3c9528e/e2e/$NEW_GUID2.md", "", "", "$NEW_GUID2.md") | Out-Null

$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:P3"))

Write-Output "handback status report regenerated"
